# Phase 3 RAD Non-UI Test Cases and Data
# Adds additional TaxType test-case rows under the existing
# "Existing Liability w/Notice Number" and "New Tax Return Amount Due"
# PaymentType blocks (columns C:E, starting at row 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Existing Liability w/Notice Number", "Admissions and Amusement Tax"),
    @("Existing Liability w/Notice Number", "Estate Tax"),
    @("Existing Liability w/Notice Number", "Motor Fuel Tax"),
    @("Existing Liability w/Notice Number", "Slots License Fee"),
    @("Existing Liability w/Notice Number", "Tobacco Tax"),
    @("Existing Liability w/Notice Number", "Transportation Network Services"),
    @("Existing Liability w/Notice Number", "Unclaimed Property"),
    @("Existing Liability w/Notice Number", "IFTA Tax"),
    @("New Tax Return Amount Due", "Admissions and Amusement Tax"),
    @("New Tax Return Amount Due", "Alcohol Tax"),
    @("New Tax Return Amount Due", "Bay Restoration Fee"),
    @("New Tax Return Amount Due", "Corporate Income Tax"),
    @("New Tax Return Amount Due", "Estate Tax"),
    @("New Tax Return Amount Due", "Motor Fuel Tax"),
    @("New Tax Return Amount Due", "Sales and Use Tax"),
    @("New Tax Return Amount Due", "Slots License Fee"),
    @("New Tax Return Amount Due", "Tire Recycling Fee"),
    @("New Tax Return Amount Due", "Tobacco Tax"),
    @("New Tax Return Amount Due", "Transportation Network Services"),
    @("New Tax Return Amount Due", "Unclaimed Property"),
    @("New Tax Return Amount Due", "Withholding Tax")
)

$row = 14
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 3).Value = "Y"
    $ws.Cells.Item($row, 4).Value = $entry[0]
    $ws.Cells.Item($row, 5).Value = $entry[1]
    $row = $row + 1
}

$ws.Range("C14:E34").Select()
